$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = 6346
$ws.Range("C24").Value = 999
$ws.Range("D24").Value = 5942061
$ws.Range("E24").Value = 936.3474629687993
$ws.Range("F24").Value = 8.182748039549947
$ws.Range("G24").Value = 3.523316062176174
$ws.Range("H24").Value = 25.87586843289331
